$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 43132
$ws.Range("A8").NumberFormat = "d-mmm"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "R "

$ws.Range("A9").Value = 43136
$ws.Range("A9").NumberFormat = "d-mmm"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "starting power and sample size calcualtions "

$ws.Range("B14").Select()
